# Update STATUS (column F) values for several activities.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Criação do repositório): Em andamento -> Concluído
$ws.Range("F3").Value = "Concluído"
# Row 4 (Check in do projeto MVC inicial no repositório): Pendente -> Concluído
$ws.Range("F4").Value = "Concluído"
# Row 5 (Check in do projeto Web Service inicial no repositório): Pendente -> Concluído
$ws.Range("F5").Value = "Concluído"
# Row 18 (Tela de cadastro/alteração de veículo): Pendente -> Em andamento
$ws.Range("F18").Value = "Em andamento"
# Row 22 (Roteiro CRUD): Pendente -> Em andamento
$ws.Range("F22").Value = "Em andamento"

# Move the active selection to the last row (A23:F23), matching the
# author's cursor position when the workbook was saved.
$ws.Range("A23:F23").Select()
